$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.119.83'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '1.826.72'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.67'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4590'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +7.70%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3738'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07334'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8607'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.01'
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = '1.825.51'
$ws.Range("E12").Value = '  +0.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.695'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.05'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +6.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.351'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07074'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008852'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.02'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '27.120.17'
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.194'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.24%  '
$ws.Range("E23").Value = '  +1.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.999'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.75'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.218'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.50'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.270'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.64'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08882'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.7733'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.196'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.972'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.469'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.40%  '
$ws.Range("E35").Value = '  -0.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.105'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01971'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("E38").Value = '  +0.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5393'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +7.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.198'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.884'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.88%  '
$ws.Range("E42").Value = '  +2.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5282'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +12.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.631'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.75'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.990'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +10.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '106.18'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06516'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.681'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9254'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.44%  '
